$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "All studies" sample query text (without the Tumor/Analyte columns) that
# replaces the old sample query previously stored in B3. Row 3 (SamplesTab)
# keeps a sample query, just an updated one.
$sampleQueryText = "SELECT`n    DISTINCT (smp.sample_id) AS ""Sample ID"",`n    sp.participant_id AS ""Participant ID"", `n    s.study_name AS ""Study Name"",`n    s.phs_accession AS Accession`nFROM `n    df_participant sp`nJOIN `n    df_study s ON sp.""study.phs_accession"" = s.phs_accession`nJOIN `n    df_sample smp ON smp.""participant.study_participant_id"" = sp.study_participant_id`nJOIN`n    df_diagnosis d ON d.""participant.study_participant_id"" = sp.study_participant_id`nJOIN`n    df_program p ON p.program_acronym = s.""program.program_acronym""`nJOIN`n    df_file f1 ON f1.""sample.sample_id"" = smp.sample_id`nJOIN`n    df_genomic_info gi ON gi.""file.file_id"" = f1.file_id`nWHERE `n    s.phs_accession = 'phs001437' AND gi.library_source = 'Transcriptomic'`nORDER BY `n    smp.sample_id ASC`nLIMIT 100;"

# Row 3 (SamplesTab) gets the refreshed sample query; row 4 (FilesTab) keeps
# holding the (unchanged) file query, just re-assigned to refresh the shared
# string pool ordering.
$fileQueryText = $ws.Range("B4").Value()
$ws.Range("B3").Value = $sampleQueryText
$ws.Range("B4").Value = $fileQueryText

# Drop the TSV/Web filename references from rows 3 and 4 (columns D and E).
$ws.Range("D3:E4").ClearContents()

# Update the saved view/selection to match the new focus (top-left A3, selection C3).
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("C3").Select()
